$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluation form")
$ws.Activate()

# Fill in the evaluation scores for each rubric row.
$ws.Range("E20").Value = 10
$ws.Range("E21").Value = 5
$ws.Range("E22").Value = 5
$ws.Range("E23").Value = 10
$ws.Range("E24").Value = 7
$ws.Range("E25").Value = 5
$ws.Range("E26").Value = 5
$ws.Range("E27").Value = 5

$ws.Range("E30").Value = 10
$ws.Range("E31").Value = 5

$ws.Range("E34").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("E36").Value = 0

# Force the overall-grade formula (which depends on the cells above) to
# recompute and refresh its cached value.
$ws.Range("E17").Formula = $ws.Range("E17").Formula

# Leave the selection where the evaluator last left off.
$ws.Range("E25").Select()
